$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new data rows (39 and 40) following the existing data pattern.
$ws.Range("A39").Value = 0
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 0.085427
$ws.Range("D39").Value = -0.1009078344957044
$ws.Range("E39").Value = "query"

$ws.Range("A40").Value = 0
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 0.170854
$ws.Range("D40").Value = -0.1825362898147634
$ws.Range("E40").Value = "query"
